$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the HTTP results row (row 7) with the new measured values.
$ws.Range("B7").Value = 11.01
$ws.Range("C7").Value = 3.13
$ws.Range("D7").Value = 735.13
$ws.Range("E7").Value = 264.9
$ws.Range("F7").Value = 2122.27
$ws.Range("G7").Value = 500.15
$ws.Range("H7").Value = 5662.31
$ws.Range("I7").Value = 1622.94
$ws.Range("J7").Value = 3
$ws.Range("K7").Value = 1.02
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1

# The row no longer carries the heavy bordered/boxed formatting - reset to
# the workbook's default "Normal" style (drops the explicit borders).
$ws.Range("B7:M7").Style = "Normal"

# Update the view: slightly higher zoom and the active selection now spans
# the whole HTTP data row instead of just the last cell.
$ws.Application.ActiveWindow.Zoom = 116
$ws.Range("B7:M7").Select() | Out-Null
